# "Generate Report for Handoff"
# The handoff report generation run updated the localization-status
# workbook: the "Priority" column now shows "ht" for the rows that were
# handed off, and the "Latest Handoff/Handback" timestamps moved forward
# a few seconds once the fresh report was produced.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows corresponding to the files that were (re)handed off in this run.
# Row 10 (7f129cd3-...) already has its full handoff/handback data and is
# left untouched.
$rows = 8, 9, 11, 12, 13, 14

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $overview.Cells.Item($r, 7).Value = "2016-08-19 16:19:52"

    # zh-cn sheet: Priority column E, and Latest Handoff Datetime column H
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-19 16:19:48"

    # de-de sheet: Priority column E, and Latest Handoff Datetime column H
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-19 16:19:52"
}
